$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "posicao" (position) column so concept position can be imported/exported.
$ws.Range("J1").Value = "posicao"
$ws.Range("J2").Value = 2
$ws.Range("J3").Value = 1
$ws.Range("J4").Value = 1

$ws.Rows.Item(1).RowHeight = 16.5
$ws.Rows.Item(2).RowHeight = 16.5
$ws.Rows.Item(3).RowHeight = 16.5
$ws.Rows.Item(4).RowHeight = 16.5

$ws.PageSetup.PaperSize = 9

$ws.Range("C6").Select()
